$d = $word.ActiveDocument

# Change 1: "This is a Markdown File" -> "This is a Markdown file"
$d.Content.Find.Execute("This is a Markdown File", $true, $true, $false, $false, $false,
                         $true, 1, $false, "This is a Markdown file", 2)

# Change 2a: Trim the trailing sentence about embedding an R code chunk from
# the "When you click the Knit button..." paragraph.
$d.Content.Find.Execute(" You can embed an R code chunk like this:", $true, $true, $false, $false, $false,
                         $true, 1, $false, "", 2)

# Change 2b: Remove the entire following paragraph (starting with
# "Note that the" and ending with "...generated the plot.") along with its
# paragraph mark.
$rng = $d.Content
$rng.Find.Execute("Note that the", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$paraStart = $rng.Paragraphs(1).Range.Start

$rng2 = $d.Content
$rng2.Find.Execute("generated the plot.", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$paraEnd = $rng2.Paragraphs(1).Range.End

$deleteRange = $d.Range($paraStart, $paraEnd)
$deleteRange.Delete()
